$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching formatting of the existing
# header row (e.g. G1: bold font, border, centered alignment)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the corresponding data value in H2
$ws.Range("H2").Value = 0
